$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for ccd11f7e-... (row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 13:13:38"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 13:13:33"
$wsZhCn.Range("K2").Value = "2016-09-01 13:14:14"

# "de-de" sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-01 13:13:38"
$wsDeDe.Range("K2").Value = "2016-09-01 13:14:24"
